# Fix uncertainty data in LCIs  (lci-biogas.xlsx, "Biogas" sheet)
#
# 1) Column M ("negative") holds a flag cell per exchange row (14-43).
#    Every row currently stores a stray literal 0 even though the column
#    is meant to be blank unless the exchange amount is actually negative.
#    Clear the bogus 0s for all rows except the one that is genuinely
#    negative (row 43).
# 2) Row 43 (exchange amount in B43 is negative) needs:
#      - H43 turned into a live formula driven off B43 instead of a
#        hard-coded literal, and
#      - M43 flagged TRUE (boolean) to mark the amount as negative.
# 3) Leave the cursor/selection where the editor ended up (H44, scrolled
#    down near row 28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biogas")

# --- 1) Drop the placeholder 0s in M14:M42 -------------------------------
$rows = 14..42
foreach ($r in $rows) {
    $ws.Range("M$r").ClearContents()
}

# --- 2) Row 43: live formula + boolean flag ------------------------------
$ws.Range("H43").Formula = "=LN(B43*-1)"
$ws.Range("M43").Value = $true

# --- 3) Restore on-screen selection/scroll position ----------------------
$ws.Range("H44").Select()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
